$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (final test score) values for specific rows, per commit "Update OOP sample final tests"
$values = @{
    4  = 200
    6  = 800
    7  = 800
    8  = 800
    9  = 800
    10 = 300
    11 = 450
    15 = 800
    16 = 700
    18 = 500
    19 = 200
    20 = 750
    21 = 300
    22 = 400
    23 = 800
    24 = 750
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
